$d = $word.ActiveDocument

$replacements = @(
    @{Old = "96×40="; New = "42×48="},
    @{Old = "90×43="; New = "22×47="},
    @{Old = "34×79="; New = "33×48="},
    @{Old = "34×82="; New = "32×91="},
    @{Old = "83×64="; New = "94×24="},
    @{Old = "38×64="; New = "45×19="},
    @{Old = "39×64="; New = "99×65="},
    @{Old = "99×28="; New = "44×43="},
    @{Old = "13×63="; New = "21×87="},
    @{Old = "13×12="; New = "12×95="},
    @{Old = "67×32="; New = "49×59="},
    @{Old = "30×25="; New = "67×49="},
    @{Old = "66×74="; New = "66×14="},
    @{Old = "30×17="; New = "50×31="},
    @{Old = "95×21="; New = "99×84="},
    @{Old = "71×14="; New = "26×52="},
    @{Old = "60×86="; New = "29×66="},
    @{Old = "16×91="; New = "49×96="},
    @{Old = "80×29="; New = "30×18="},
    @{Old = "87×93="; New = "22×35="},
    @{Old = "46×67="; New = "38×93="},
    @{Old = "84×67="; New = "81×83="},
    @{Old = "55×46="; New = "39×38="},
    @{Old = "86×81="; New = "59×60="},
    @{Old = "51×99="; New = "70×40="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
